$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 changes from PREPROD to QA
$ws.Range("A5").Value = "QA"

# New row 12 (PREPROD group gains one more claim number)
$ws.Range("A12").Value = "PREPROD"

# Update cell values in the same order the strings were (re)typed so the
# resulting shared-strings table lines up with the saved workbook.
$ws.Range("B6").Value = "'0420172010219    "
$ws.Range("B8").Value = "'0420172010221"
$ws.Range("B9").Value = "'0420172010222"
$ws.Range("B10").Value = "'1220170301466"
$ws.Range("B11").Value = "'1120170200969"
$ws.Range("B12").Value = "'1220170301467"
$ws.Range("B7").Value = "'0420172010228"
$ws.Range("B2").Value = "'0420194406895   "
$ws.Range("B3").Value = "'0420194406896"
$ws.Range("B4").Value = "'1120170200973"
$ws.Range("B5").Value = "'1220194200691"

# Update selection to match the recorded UI state
$ws.Range("A6:XFD6").Select()
